# Scheduled-runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for a batch of
# leves across the Golem_Profits sheets, per the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 199.5
$ws.Range("J2").Value = 200
$ws.Range("L2").Value = 200
$ws.Range("N2").Value = -426

$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 102.10345
$ws.Range("I33").Value = 75.96296
$ws.Range("K33").Value = 75.96296
$ws.Range("M33").Value = 153.03704

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H88").Value = 3127.4614
$ws.Range("J88").Value = 3218.818
$ws.Range("L88").Value = 3218.818
$ws.Range("N88").Value = -4030.818

$ws.Range("H91").Value = 3127.4614
$ws.Range("J91").Value = 3218.818
$ws.Range("L91").Value = 3218.818
$ws.Range("N91").Value = -6026.818

$ws.Range("H106").Value = 1866.3334
$ws.Range("I106").Value = 1866.3334
$ws.Range("K106").Value = 1866.3334
$ws.Range("M106").Value = -1235.3334

$ws.Range("H135").Value = 868.75
$ws.Range("I135").Value = 868.75
$ws.Range("K135").Value = 7818.75
$ws.Range("M135").Value = -5283.75

$ws.Range("H138").Value = 2247.0667
$ws.Range("I138").Value = 728.3333
$ws.Range("J138").Value = 3259.5557
$ws.Range("K138").Value = 2184.9999
$ws.Range("L138").Value = 9778.667099999999
$ws.Range("M138").Value = 2955.0001
$ws.Range("N138").Value = -20058.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H22").Value = 1072
$ws.Range("I22").Value = 1110.5
$ws.Range("K22").Value = 1110.5
$ws.Range("M22").Value = -811.5

$ws.Range("H23").Value = 9999
$ws.Range("J23").Value = 9999
$ws.Range("L23").Value = 9999
$ws.Range("N23").Value = -10517

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

$ws.Range("H124").Value = 35688.8
$ws.Range("J124").Value = 35688.8
$ws.Range("L124").Value = 35688.8
$ws.Range("N124").Value = -45508.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 20050
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 40000
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 40000
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -40226

$ws.Range("H26").Value = 23964
$ws.Range("I26").Value = 23964
$ws.Range("K26").Value = 23964
$ws.Range("M26").Value = -23672

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H86").Value = 6944.778
$ws.Range("I86").Value = 5917.5
$ws.Range("K86").Value = 5917.5
$ws.Range("M86").Value = -4794.5

$ws.Range("H89").Value = 6944.778
$ws.Range("I89").Value = 5917.5
$ws.Range("K89").Value = 29587.5
$ws.Range("M89").Value = -23971.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2256.111
$ws.Range("J2").Value = 4135
$ws.Range("L2").Value = 4135
$ws.Range("N2").Value = -4361

$ws.Range("H7").Value = 57.7
$ws.Range("J7").Value = 53.833332
$ws.Range("L7").Value = 53.833332
$ws.Range("N7").Value = -279.833332

$ws.Range("H22").Value = 775.75
$ws.Range("I22").Value = 775.75
$ws.Range("K22").Value = 775.75
$ws.Range("M22").Value = -425.75

$ws.Range("H31").Value = 1350
$ws.Range("J31").Value = 1350
$ws.Range("L31").Value = 1350
$ws.Range("N31").Value = -1940

$ws.Range("H34").Value = 1350
$ws.Range("J34").Value = 1350
$ws.Range("L34").Value = 1350
$ws.Range("N34").Value = -1754

$ws.Range("H105").Value = 711.3333
$ws.Range("I105").Value = 653.6
$ws.Range("K105").Value = 653.6
$ws.Range("M105").Value = 1093.4

$ws.Range("H125").Value = 15749.75
$ws.Range("J125").Value = 15749.75
$ws.Range("L125").Value = 15749.75
$ws.Range("N125").Value = -20669.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23.272728
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 31
$ws.Range("K2").Value = 84
$ws.Range("L2").Value = 186
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = -412

$ws.Range("H4").Value = 2548.5366
$ws.Range("I4").Value = 1874.6428
$ws.Range("K4").Value = 5623.928400000001
$ws.Range("M4").Value = -5511.928400000001

$ws.Range("H12").Value = 171.4
$ws.Range("J12").Value = 237.42857
$ws.Range("L12").Value = 712.28571
$ws.Range("N12").Value = -1058.28571

$ws.Range("H15").Value = 36
$ws.Range("I15").Value = 32.8
$ws.Range("J15").Value = 44
$ws.Range("K15").Value = 98.39999999999999
$ws.Range("L15").Value = 132
$ws.Range("M15").Value = 41.60000000000001
$ws.Range("N15").Value = -412

$ws.Range("H41").Value = 1166
$ws.Range("I41").Value = 999
$ws.Range("J41").Value = 1249.5
$ws.Range("K41").Value = 2997
$ws.Range("L41").Value = 3748.5
$ws.Range("M41").Value = -2659
$ws.Range("N41").Value = -4424.5

$ws.Range("H44").Value = 1117.7778
$ws.Range("I44").Value = 580
$ws.Range("K44").Value = 1740
$ws.Range("M44").Value = -1342

$ws.Range("H113").Value = 199.5
$ws.Range("I113").Value = 199
$ws.Range("K113").Value = 597
$ws.Range("M113").Value = 1573

$ws.Range("H131").Value = 4974.5
$ws.Range("J131").Value = 4974.5
$ws.Range("L131").Value = 14923.5
$ws.Range("N131").Value = -25003.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5916666.5
$ws.Range("J3").Value = 6666666.5
$ws.Range("L3").Value = 6666666.5
$ws.Range("N3").Value = -6666898.5

$ws.Range("H5").Value = 998.5
$ws.Range("I5").Value = 998.5
$ws.Range("K5").Value = 998.5
$ws.Range("M5").Value = -886.5

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3668
$ws.Range("I122").Value = 3502
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10506
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8056
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H100").Value = 199.5
$ws.Range("I100").Value = 199.5
$ws.Range("K100").Value = 399
$ws.Range("M100").Value = 142

$ws.Range("H136").Value = 6109.5557
$ws.Range("I136").Value = 6109.5557
$ws.Range("K136").Value = 18328.6671
$ws.Range("M136").Value = -15778.6671
